# Integrate Case 3 (PDF Table extraction) into the code base
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Georgia) updates: refreshed daily-pull numbers ---
$ws.Range("D5").Value = 33508
$ws.Range("E5").Value = 1405
$ws.Range("F5").Value = 11857
$ws.Range("H5").Value = 35.39
$ws.Range("I5").Value = 49.61

# --- Row 6 (Michigan) updates: refreshed daily-pull numbers + date ---
$ws.Cells.Item(6, 3).NumberFormat = "@"
$ws.Cells.Item(6, 3).Value = "5/10/2020"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Range("D6").Value = 47138
$ws.Range("E6").Value = 4551
$ws.Range("F6").Value = 15084
$ws.Range("G6").Value = 1866

# --- Row 10 (Wisconsin -- Milwaukee) updates: refreshed daily-pull numbers ---
$ws.Range("D10").Value = 3981
$ws.Range("E10").Value = 217
$ws.Range("F10").Value = 1544
$ws.Range("G10").Value = 106
$ws.Range("H10").Value = 38.78
$ws.Range("I10").Value = 48.85

# --- New Row 11 (San Diego) : first successful PDF-table extraction (Case 3) ---
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(11, 1).Value = "San Diego"
$ws.Range("B11").Value = "California - San Diego"
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "5/9/2020"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Range("D11").Value = 4926
$ws.Range("E11").Value = 175
$ws.Range("F11").Value = 167
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 3.39
$ws.Range("I11").Value = 2.29
$ws.Range("J11").Value = "Success!"

# --- New Row 12 (Florida) : Case 3 extraction failed -- no PDF found ---
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)
$ws.Cells.Item(12, 1).Value = "Florida"
$ws.Range("B12").Value = "Florida"
$ws.Range("J12").Value = "An error occured. ... FileNotFoundError(2, 'No such file or directory')"
